# Natmi following Dr Hou advice
#
# The LR-pairs table (Efnb2 -> Epha4) is rebuilt: instead of limiting the
# "Target cluster" to each sending cluster's two historical partners, every
# sending cluster (ECs, FAPs, M2, sCs) is now paired with every possible
# target cluster (ECs, FAPs, M2, sCs) -- a full 4x4 cross-product (16 rows,
# A2:T17) -- with updated ligand/receptor expression statistics recomputed
# for each pairing.

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("A2").Value = "ECs"
$ws.Range("B2").Value = "Efnb2"
$ws.Range("C2").Value = "Epha4"
$ws.Range("D2").Value = "ECs"
$ws.Range("E2").Value = 3
$ws.Range("F2").Value = 1
$ws.Range("G2").Value = 28.95628266666667
$ws.Range("H2").Value = 86.868848
$ws.Range("I2").Value = 0.5491054194301004
$ws.Range("J2").Value = 0.5491054194301005
$ws.Range("K2").Value = 2
$ws.Range("L2").Value = 0.6666666666666666
$ws.Range("M2").Value = 5.685057
$ws.Range("N2").Value = 17.055171
$ws.Range("O2").Value = 0.3604606774420115
$ws.Range("P2").Value = 0.3604606774420115
$ws.Range("Q2").Value = 164.618117468112
$ws.Range("R2").Value = 1481.563057213008
$ws.Range("S2").Value = 0.1979309114748539
$ws.Range("T2").Value = 0.1979309114748539
$ws.Range("A3").Value = "ECs"
$ws.Range("B3").Value = "Efnb2"
$ws.Range("C3").Value = "Epha4"
$ws.Range("D3").Value = "FAPs"
$ws.Range("E3").Value = 3
$ws.Range("F3").Value = 1
$ws.Range("G3").Value = 28.95628266666667
$ws.Range("H3").Value = 86.868848
$ws.Range("I3").Value = 0.5491054194301004
$ws.Range("J3").Value = 0.5491054194301005
$ws.Range("K3").Value = 3
$ws.Range("L3").Value = 1
$ws.Range("M3").Value = 8.775186333333332
$ws.Range("N3").Value = 26.325559
$ws.Range("O3").Value = 0.556390131249909
$ws.Range("P3").Value = 0.5563901312499091
$ws.Range("Q3").Value = 254.0967759206702
$ws.Range("R3").Value = 2286.870983286032
$ws.Range("S3").Value = 0.3055168363867499
$ws.Range("T3").Value = 0.30551683638675
$ws.Range("A4").Value = "ECs"
$ws.Range("B4").Value = "Efnb2"
$ws.Range("C4").Value = "Epha4"
$ws.Range("D4").Value = "M2"
$ws.Range("E4").Value = 3
$ws.Range("F4").Value = 1
$ws.Range("G4").Value = 28.95628266666667
$ws.Range("H4").Value = 86.868848
$ws.Range("I4").Value = 0.5491054194301004
$ws.Range("J4").Value = 0.5491054194301005
$ws.Range("K4").Value = 1
$ws.Range("L4").Value = 0.3333333333333333
$ws.Range("M4").Value = 0.04169666666666667
$ws.Range("N4").Value = 0.12509
$ws.Range("O4").Value = 0.002643774497553922
$ws.Range("P4").Value = 0.002643774497553922
$ws.Range("Q4").Value = 1.207380466257778
$ws.Range("R4").Value = 10.86642419632
$ws.Range("S4").Value = 0.001451710904357949
$ws.Range("T4").Value = 0.00145171090435795
$ws.Range("A5").Value = "ECs"
$ws.Range("B5").Value = "Efnb2"
$ws.Range("C5").Value = "Epha4"
$ws.Range("D5").Value = "sCs"
$ws.Range("E5").Value = 3
$ws.Range("F5").Value = 1
$ws.Range("G5").Value = 28.95628266666667
$ws.Range("H5").Value = 86.868848
$ws.Range("I5").Value = 0.5491054194301004
$ws.Range("J5").Value = 0.5491054194301005
$ws.Range("K5").Value = 3
$ws.Range("L5").Value = 1
$ws.Range("M5").Value = 1.269702666666667
$ws.Range("N5").Value = 3.809108
$ws.Range("O5").Value = 0.08050541681052542
$ws.Range("P5").Value = 0.08050541681052542
$ws.Range("Q5").Value = 36.76586931862045
$ws.Range("R5").Value = 330.892823867584
$ws.Range("S5").Value = 0.04420596066413861
$ws.Range("T5").Value = 0.04420596066413862
$ws.Range("A6").Value = "FAPs"
$ws.Range("B6").Value = "Efnb2"
$ws.Range("C6").Value = "Epha4"
$ws.Range("D6").Value = "ECs"
$ws.Range("E6").Value = 3
$ws.Range("F6").Value = 1
$ws.Range("G6").Value = 12.691493
$ws.Range("H6").Value = 38.074479
$ws.Range("I6").Value = 0.2406720388519202
$ws.Range("J6").Value = 0.2406720388519202
$ws.Range("K6").Value = 2
$ws.Range("L6").Value = 0.6666666666666666
$ws.Range("M6").Value = 5.685057
$ws.Range("N6").Value = 17.055171
$ws.Range("O6").Value = 0.3604606774420115
$ws.Range("P6").Value = 0.3604606774420115
$ws.Range("Q6").Value = 72.151861120101
$ws.Range("R6").Value = 649.366750080909
$ws.Range("S6").Value = 0.08675280616591327
$ws.Range("T6").Value = 0.08675280616591327
$ws.Range("A7").Value = "FAPs"
$ws.Range("B7").Value = "Efnb2"
$ws.Range("C7").Value = "Epha4"
$ws.Range("D7").Value = "FAPs"
$ws.Range("E7").Value = 3
$ws.Range("F7").Value = 1
$ws.Range("G7").Value = 12.691493
$ws.Range("H7").Value = 38.074479
$ws.Range("I7").Value = 0.2406720388519202
$ws.Range("J7").Value = 0.2406720388519202
$ws.Range("K7").Value = 3
$ws.Range("L7").Value = 1
$ws.Range("M7").Value = 8.775186333333332
$ws.Range("N7").Value = 26.325559
$ws.Range("O7").Value = 0.556390131249909
$ws.Range("P7").Value = 0.5563901312499091
$ws.Range("Q7").Value = 111.3702159231956
$ws.Range("R7").Value = 1002.331943308761
$ws.Range("S7").Value = 0.1339075472850031
$ws.Range("T7").Value = 0.1339075472850031
$ws.Range("A8").Value = "FAPs"
$ws.Range("B8").Value = "Efnb2"
$ws.Range("C8").Value = "Epha4"
$ws.Range("D8").Value = "M2"
$ws.Range("E8").Value = 3
$ws.Range("F8").Value = 1
$ws.Range("G8").Value = 12.691493
$ws.Range("H8").Value = 38.074479
$ws.Range("I8").Value = 0.2406720388519202
$ws.Range("J8").Value = 0.2406720388519202
$ws.Range("K8").Value = 1
$ws.Range("L8").Value = 0.3333333333333333
$ws.Range("M8").Value = 0.04169666666666667
$ws.Range("N8").Value = 0.12509
$ws.Range("O8").Value = 0.002643774497553922
$ws.Range("P8").Value = 0.002643774497553922
$ws.Range("Q8").Value = 0.5291929531233334
$ws.Range("R8").Value = 4.76273657811
$ws.Range("S8").Value = 0.0006362825985910132
$ws.Range("T8").Value = 0.0006362825985910133
$ws.Range("A9").Value = "FAPs"
$ws.Range("B9").Value = "Efnb2"
$ws.Range("C9").Value = "Epha4"
$ws.Range("D9").Value = "sCs"
$ws.Range("E9").Value = 3
$ws.Range("F9").Value = 1
$ws.Range("G9").Value = 12.691493
$ws.Range("H9").Value = 38.074479
$ws.Range("I9").Value = 0.2406720388519202
$ws.Range("J9").Value = 0.2406720388519202
$ws.Range("K9").Value = 3
$ws.Range("L9").Value = 1
$ws.Range("M9").Value = 1.269702666666667
$ws.Range("N9").Value = 3.809108
$ws.Range("O9").Value = 0.08050541681052542
$ws.Range("P9").Value = 0.08050541681052542
$ws.Range("Q9").Value = 16.11442250608133
$ws.Range("R9").Value = 145.029802554732
$ws.Range("S9").Value = 0.0193754028024128
$ws.Range("T9").Value = 0.0193754028024128
$ws.Range("A10").Value = "M2"
$ws.Range("B10").Value = "Efnb2"
$ws.Range("C10").Value = "Epha4"
$ws.Range("D10").Value = "ECs"
$ws.Range("E10").Value = 3
$ws.Range("F10").Value = 1
$ws.Range("G10").Value = 0.4888703333333334
$ws.Range("H10").Value = 1.466611
$ws.Range("I10").Value = 0.009270573592685367
$ws.Range("J10").Value = 0.009270573592685367
$ws.Range("K10").Value = 2
$ws.Range("L10").Value = 0.6666666666666666
$ws.Range("M10").Value = 5.685057
$ws.Range("N10").Value = 17.055171
$ws.Range("O10").Value = 0.3604606774420115
$ws.Range("P10").Value = 0.3604606774420115
$ws.Range("Q10").Value = 2.779255710609
$ws.Range("R10").Value = 25.013301395481
$ws.Range("S10").Value = 0.00334167723749539
$ws.Range("T10").Value = 0.00334167723749539
$ws.Range("A11").Value = "M2"
$ws.Range("B11").Value = "Efnb2"
$ws.Range("C11").Value = "Epha4"
$ws.Range("D11").Value = "FAPs"
$ws.Range("E11").Value = 3
$ws.Range("F11").Value = 1
$ws.Range("G11").Value = 0.4888703333333334
$ws.Range("H11").Value = 1.466611
$ws.Range("I11").Value = 0.009270573592685367
$ws.Range("J11").Value = 0.009270573592685367
$ws.Range("K11").Value = 3
$ws.Range("L11").Value = 1
$ws.Range("M11").Value = 8.775186333333332
$ws.Range("N11").Value = 26.325559
$ws.Range("O11").Value = 0.556390131249909
$ws.Range("P11").Value = 0.5563901312499091
$ws.Range("Q11").Value = 4.289928267838778
$ws.Range("R11").Value = 38.609354410549
$ws.Range("S11").Value = 0.005158055657996151
$ws.Range("T11").Value = 0.005158055657996152
$ws.Range("A12").Value = "M2"
$ws.Range("B12").Value = "Efnb2"
$ws.Range("C12").Value = "Epha4"
$ws.Range("D12").Value = "M2"
$ws.Range("E12").Value = 3
$ws.Range("F12").Value = 1
$ws.Range("G12").Value = 0.4888703333333334
$ws.Range("H12").Value = 1.466611
$ws.Range("I12").Value = 0.009270573592685367
$ws.Range("J12").Value = 0.009270573592685367
$ws.Range("K12").Value = 1
$ws.Range("L12").Value = 0.3333333333333333
$ws.Range("M12").Value = 0.04169666666666667
$ws.Range("N12").Value = 0.12509
$ws.Range("O12").Value = 0.002643774497553922
$ws.Range("P12").Value = 0.002643774497553922
$ws.Range("Q12").Value = 0.02038426333222222
$ws.Range("R12").Value = 0.18345836999
$ws.Range("S12").Value = [double]"2.450930604203841e-05"
$ws.Range("T12").Value = [double]"2.450930604203842e-05"
$ws.Range("A13").Value = "M2"
$ws.Range("B13").Value = "Efnb2"
$ws.Range("C13").Value = "Epha4"
$ws.Range("D13").Value = "sCs"
$ws.Range("E13").Value = 3
$ws.Range("F13").Value = 1
$ws.Range("G13").Value = 0.4888703333333334
$ws.Range("H13").Value = 1.466611
$ws.Range("I13").Value = 0.009270573592685367
$ws.Range("J13").Value = 0.009270573592685367
$ws.Range("K13").Value = 3
$ws.Range("L13").Value = 1
$ws.Range("M13").Value = 1.269702666666667
$ws.Range("N13").Value = 3.809108
$ws.Range("O13").Value = 0.08050541681052542
$ws.Range("P13").Value = 0.08050541681052542
$ws.Range("Q13").Value = 0.6207199658875555
$ws.Range("R13").Value = 5.586479692988001
$ws.Range("S13").Value = 0.0007463313911517856
$ws.Range("T13").Value = 0.0007463313911517856
$ws.Range("A14").Value = "sCs"
$ws.Range("B14").Value = "Efnb2"
$ws.Range("C14").Value = "Epha4"
$ws.Range("D14").Value = "ECs"
$ws.Range("E14").Value = 3
$ws.Range("F14").Value = 1
$ws.Range("G14").Value = 10.59691233333333
$ws.Range("H14").Value = 31.790737
$ws.Range("I14").Value = 0.200951968125294
$ws.Range("J14").Value = 0.200951968125294
$ws.Range("K14").Value = 2
$ws.Range("L14").Value = 0.6666666666666666
$ws.Range("M14").Value = 5.685057
$ws.Range("N14").Value = 17.055171
$ws.Range("O14").Value = 0.3604606774420115
$ws.Range("P14").Value = 0.3604606774420115
$ws.Range("Q14").Value = 60.24405063900301
$ws.Range("R14").Value = 542.1964557510271
$ws.Range("S14").Value = 0.072435282563749
$ws.Range("T14").Value = 0.072435282563749
$ws.Range("A15").Value = "sCs"
$ws.Range("B15").Value = "Efnb2"
$ws.Range("C15").Value = "Epha4"
$ws.Range("D15").Value = "FAPs"
$ws.Range("E15").Value = 3
$ws.Range("F15").Value = 1
$ws.Range("G15").Value = 10.59691233333333
$ws.Range("H15").Value = 31.790737
$ws.Range("I15").Value = 0.200951968125294
$ws.Range("J15").Value = 0.200951968125294
$ws.Range("K15").Value = 3
$ws.Range("L15").Value = 1
$ws.Range("M15").Value = 8.775186333333332
$ws.Range("N15").Value = 26.325559
$ws.Range("O15").Value = 0.556390131249909
$ws.Range("P15").Value = 0.5563901312499091
$ws.Range("Q15").Value = 92.9898802829981
$ws.Range("R15").Value = 836.9089225469829
$ws.Range("S15").Value = 0.1118076919201599
$ws.Range("T15").Value = 0.1118076919201599
$ws.Range("A16").Value = "sCs"
$ws.Range("B16").Value = "Efnb2"
$ws.Range("C16").Value = "Epha4"
$ws.Range("D16").Value = "M2"
$ws.Range("E16").Value = 3
$ws.Range("F16").Value = 1
$ws.Range("G16").Value = 10.59691233333333
$ws.Range("H16").Value = 31.790737
$ws.Range("I16").Value = 0.200951968125294
$ws.Range("J16").Value = 0.200951968125294
$ws.Range("K16").Value = 1
$ws.Range("L16").Value = 0.3333333333333333
$ws.Range("M16").Value = 0.04169666666666667
$ws.Range("N16").Value = 0.12509
$ws.Range("O16").Value = 0.002643774497553922
$ws.Range("P16").Value = 0.002643774497553922
$ws.Range("Q16").Value = 0.4418559212588889
$ws.Range("R16").Value = 3.97670329133
$ws.Range("S16").Value = 0.0005312716885629209
$ws.Range("T16").Value = 0.000531271688562921
$ws.Range("A17").Value = "sCs"
$ws.Range("B17").Value = "Efnb2"
$ws.Range("C17").Value = "Epha4"
$ws.Range("D17").Value = "sCs"
$ws.Range("E17").Value = 3
$ws.Range("F17").Value = 1
$ws.Range("G17").Value = 10.59691233333333
$ws.Range("H17").Value = 31.790737
$ws.Range("I17").Value = 0.200951968125294
$ws.Range("J17").Value = 0.200951968125294
$ws.Range("K17").Value = 3
$ws.Range("L17").Value = 1
$ws.Range("M17").Value = 1.269702666666667
$ws.Range("N17").Value = 3.809108
$ws.Range("O17").Value = 0.08050541681052542
$ws.Range("P17").Value = 0.08050541681052542
$ws.Range("Q17").Value = 13.45492784806622
$ws.Range("R17").Value = 121.094350632596
$ws.Range("S17").Value = 0.01617772195282222
$ws.Range("T17").Value = 0.01617772195282222
